$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "25.950.89"
$ws.Range("E2").Value2 = "  -0.26%  "

$ws.Range("D3").Value2 = "1.644.19"
$ws.Range("E3").Value2 = "  +0.24%  "

$ws.Range("E4").Value2 = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "215.68"
$ws.Range("E5").Value2 = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5058"
$ws.Range("E6").Value2 = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.005"
$ws.Range("E7").Value2 = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2582"
$ws.Range("E8").Value2 = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06413"
$ws.Range("E9").Value2 = "  -0.51%  "

$ws.Range("E10").Value2 = "  +0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07784"
$ws.Range("E11").Value2 = "  +0.96%  "

$ws.Range("D12").Value2 = "1.659.20"
$ws.Range("E12").Value2 = "  +1.05%  "

$ws.Range("E13").Value2 = "  +1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.5448"
$ws.Range("E14").Value2 = "  +0.03%  "

$ws.Range("D15").Value2 = "0.0₅7887"
$ws.Range("E15").Value2 = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "65.01"

$ws.Range("D17").Value2 = "25.991.17"
$ws.Range("E17").Value2 = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "198.62"
$ws.Range("E19").Value2 = "  -2.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "4.415"
$ws.Range("E20").Value2 = "  +3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.989"
$ws.Range("E21").Value2 = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.007"
$ws.Range("E22").Value2 = "  +0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.007"
$ws.Range("E23").Value2 = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "1.871"
$ws.Range("E24").Value2 = "  -4.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "141.11"
$ws.Range("E25").Value2 = "  -0.46%  "

$ws.Range("E26").Value2 = "  -0.24%  "

$ws.Range("E27").Value2 = "  +2.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "15.76"
$ws.Range("E28").Value2 = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.245"
$ws.Range("E29").Value2 = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.05005"
$ws.Range("E30").Value2 = "  -0.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.274"
$ws.Range("E31").Value2 = "  +0.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.207"
$ws.Range("E32").Value2 = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.535"
$ws.Range("E33").Value2 = "  -0.25%  "

$ws.Range("E34").Value2 = "  +1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.8963"
$ws.Range("E35").Value2 = "  +0.72%  "

$ws.Range("E36").Value2 = "  -0.77%  "

$ws.Range("D37").Value2 = "1.145.68"
$ws.Range("E37").Value2 = "  -0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.5561"
$ws.Range("E38").Value2 = "  -1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01566"
$ws.Range("E39").Value2 = "  -0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.006"
$ws.Range("E40").Value2 = "  -0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "5.694"
$ws.Range("E41").Value2 = "  +0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.8248"
$ws.Range("E42").Value2 = "  +2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "99.99"
$ws.Range("E43").Value2 = "  +0.28%  "

$ws.Range("E44").Value2 = "  +6.55%  "

$ws.Range("D45").Value2 = "1.782.16"
$ws.Range("E45").Value2 = "  +0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.4530"
$ws.Range("E46").Value2 = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "55.51"
$ws.Range("E47").Value2 = "  +0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.005"
$ws.Range("E48").Value2 = "  -0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.05062"
$ws.Range("E49").Value2 = "  +0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.007"
$ws.Range("E50").Value2 = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.09554"
$ws.Range("E51").Value2 = "  +2.69%  "
